# Fix Training Data Issue (#48)
# Data in column BF ("Date") was off by one day because of how the NBA
# stats site displayed dates. The original text "6-27-2013-14" is
# replaced with the corrected ISO-style date "2014-06-27" for every data
# row (rows 2-31 -- the header in row 1 is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 31
$col      = "BF"
$newDate  = "2014-06-27"

$rng = $ws.Range("$col$firstRow`:$col$lastRow")

# Force a text number format first so that Excel's automatic date
# recognition does not silently convert the literal string we are about
# to write into a date serial number (the source data must stay plain
# text, exactly like the other string columns in this sheet).
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("$col$r").Value2 = $newDate
}

# Restore the default style on the range so the cells keep behaving like
# ordinary unstyled text cells (matches the rest of the sheet).
$rng.Style = "Normal"
